# Update the statsmodels OLS summary text embedded in each sheet's B2 cell.
# Every sheet ("18".."5") carries a long pre-formatted text report in B2 with
# a "Date:" line ("Wed, 01 Jan 2020") and a "Time:" line ("23:18:45"/"23:18:46").
# Both get bumped to the re-run's timestamp: Thu, 02 Jan 2020 / 20:48:37.

$wb = $excel.ActiveWorkbook

$oldDate = "Wed, 01 Jan 2020"
$newDate = "Thu, 02 Jan 2020"
$oldTimes = @("23:18:45", "23:18:46")
$newTime = "20:48:37"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -ne $null -and $text.Contains("Date:")) {
        $updated = $text.Replace($oldDate, $newDate)
        foreach ($ot in $oldTimes) {
            $updated = $updated.Replace($ot, $newTime)
        }

        if ($updated -ne $text) {
            $cell.Value2 = $updated
        }
    }
}

Write-Output "Updated Date/Time stamps in OLS summary reports across all sheets."
